$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Name" row value changes to the technical/code name
$ws.Range("B4").Value = "ValueSetCBO"

# Insert the previous Name value into the existing "Title" row
$ws.Range("B5").Value = "Classificação Brasileira de Ocupações - CBO"
